$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether the value must be forced to
# Text (column D sometimes holds numeric-looking strings like "6.23" that Excel
# would otherwise auto-convert to a real number instead of leaving as text).
$cellUpdates = @(
    @{ Cell = "D2"; Value = "58.329.60"; ForceText = $false }
    @{ Cell = "E2"; Value = "  +2.03%  "; ForceText = $false }
    @{ Cell = "D3"; Value = "3.155.97"; ForceText = $false }
    @{ Cell = "E3"; Value = "  +2.72%  "; ForceText = $false }
    @{ Cell = "E4"; Value = "  -0.03%  "; ForceText = $false }
    @{ Cell = "D5"; Value = "537.28"; ForceText = $true }
    @{ Cell = "E5"; Value = "  +2.93%  "; ForceText = $false }
    @{ Cell = "D6"; Value = "139.94"; ForceText = $true }
    @{ Cell = "E6"; Value = "  +3.28%  "; ForceText = $false }
    @{ Cell = "E7"; Value = "  +0.00%  "; ForceText = $false }
    @{ Cell = "D8"; Value = "0.516"; ForceText = $true }
    @{ Cell = "E8"; Value = "  +10.03%  "; ForceText = $false }
    @{ Cell = "D9"; Value = "7.34"; ForceText = $true }
    @{ Cell = "E9"; Value = "  +0.57%  "; ForceText = $false }
    @{ Cell = "D10"; Value = "0.110"; ForceText = $true }
    @{ Cell = "E10"; Value = "  +3.24%  "; ForceText = $false }
    @{ Cell = "D11"; Value = "0.423"; ForceText = $true }
    @{ Cell = "E11"; Value = "  +5.45%  "; ForceText = $false }
    @{ Cell = "E12"; Value = "  +2.97%  "; ForceText = $false }
    @{ Cell = "D13"; Value = "3.693.51"; ForceText = $false }
    @{ Cell = "E13"; Value = "  +2.76%  "; ForceText = $false }
    @{ Cell = "D14"; Value = "26.11"; ForceText = $true }
    @{ Cell = "E14"; Value = "  +3.77%  "; ForceText = $false }
    @{ Cell = "E15"; Value = "  +5.65%  "; ForceText = $false }
    @{ Cell = "D16"; Value = "58.372.37"; ForceText = $false }
    @{ Cell = "E16"; Value = "  +2.00%  "; ForceText = $false }
    @{ Cell = "B17"; Value = "WrappedEther"; ForceText = $false }
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; ForceText = $false }
    @{ Cell = "D17"; Value = "3.144.53"; ForceText = $false }
    @{ Cell = "E17"; Value = "  +2.19%  "; ForceText = $false }
    @{ Cell = "B18"; Value = "Polkadot"; ForceText = $false }
    @{ Cell = "C18"; Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"; ForceText = $false }
    @{ Cell = "D18"; Value = "6.23"; ForceText = $true }
    @{ Cell = "E18"; Value = "  +6.27%  "; ForceText = $false }
    @{ Cell = "D19"; Value = "13.02"; ForceText = $true }
    @{ Cell = "E19"; Value = "  +4.79%  "; ForceText = $false }
    @{ Cell = "D20"; Value = "8.20"; ForceText = $true }
    @{ Cell = "E20"; Value = "  +4.90%  "; ForceText = $false }
    @{ Cell = "D21"; Value = "377.05"; ForceText = $true }
    @{ Cell = "E21"; Value = "  +8.05%  "; ForceText = $false }
    @{ Cell = "D22"; Value = "0.998"; ForceText = $true }
    @{ Cell = "E22"; Value = "  -0.12%  "; ForceText = $false }
    @{ Cell = "D23"; Value = "5.75"; ForceText = $true }
    @{ Cell = "E23"; Value = "  -0.47%  "; ForceText = $false }
    @{ Cell = "D24"; Value = "70.53"; ForceText = $true }
    @{ Cell = "E24"; Value = "  +2.25%  "; ForceText = $false }
    @{ Cell = "D25"; Value = "0.518"; ForceText = $true }
    @{ Cell = "E25"; Value = "  +4.06%  "; ForceText = $false }
    @{ Cell = "D26"; Value = "0.168"; ForceText = $true }
    @{ Cell = "E26"; Value = "  +1.58%  "; ForceText = $false }
    @{ Cell = "E27"; Value = "  +0.29%  "; ForceText = $false }
    @{ Cell = "D28"; Value = "8.13"; ForceText = $true }
    @{ Cell = "E28"; Value = "  +13.55%  "; ForceText = $false }
    @{ Cell = "D29"; Value = "0.0₃0883"; ForceText = $false }
    @{ Cell = "E29"; Value = "  +2.79%  "; ForceText = $false }
    @{ Cell = "E30"; Value = "  +2.43%  "; ForceText = $false }
    @{ Cell = "D31"; Value = "6.18"; ForceText = $true }
    @{ Cell = "E31"; Value = "  +6.78%  "; ForceText = $false }
    @{ Cell = "D32"; Value = "21.84"; ForceText = $true }
    @{ Cell = "E32"; Value = "  +4.43%  "; ForceText = $false }
    @{ Cell = "E33"; Value = "  +7.39%  "; ForceText = $false }
    @{ Cell = "E34"; Value = "  +4.48%  "; ForceText = $false }
    @{ Cell = "D35"; Value = "161.41"; ForceText = $true }
    @{ Cell = "E35"; Value = "  +1.47%  "; ForceText = $false }
    @{ Cell = "D36"; Value = "6.26"; ForceText = $true }
    @{ Cell = "E36"; Value = "  +4.80%  "; ForceText = $false }
    @{ Cell = "E37"; Value = "  +10.28%  "; ForceText = $false }
    @{ Cell = "D38"; Value = "25.54"; ForceText = $true }
    @{ Cell = "E38"; Value = "  +0.94%  "; ForceText = $false }
    @{ Cell = "E39"; Value = "  +7.92%  "; ForceText = $false }
    @{ Cell = "D40"; Value = "2.646.97"; ForceText = $false }
    @{ Cell = "E40"; Value = "  +9.79%  "; ForceText = $false }
    @{ Cell = "D41"; Value = "0.0681"; ForceText = $true }
    @{ Cell = "E41"; Value = "  +4.03%  "; ForceText = $false }
    @{ Cell = "E42"; Value = "  +5.40%  "; ForceText = $false }
    @{ Cell = "D43"; Value = "38.65"; ForceText = $true }
    @{ Cell = "E43"; Value = "  +6.00%  "; ForceText = $false }
    @{ Cell = "D44"; Value = "0.704"; ForceText = $true }
    @{ Cell = "E44"; Value = "  +1.60%  "; ForceText = $false }
    @{ Cell = "D45"; Value = "0.0277"; ForceText = $true }
    @{ Cell = "E45"; Value = "  +5.76%  "; ForceText = $false }
    @{ Cell = "E46"; Value = "  +0.00%  "; ForceText = $false }
    @{ Cell = "E47"; Value = "  +12.46%  "; ForceText = $false }
    @{ Cell = "D48"; Value = "6.24"; ForceText = $true }
    @{ Cell = "E48"; Value = "  +4.76%  "; ForceText = $false }
    @{ Cell = "D49"; Value = "0.983"; ForceText = $true }
    @{ Cell = "E49"; Value = "  +5.26%  "; ForceText = $false }
    @{ Cell = "D50"; Value = "20.31"; ForceText = $true }
    @{ Cell = "E50"; Value = "  +4.26%  "; ForceText = $false }
    @{ Cell = "D51"; Value = "0.754"; ForceText = $true }
    @{ Cell = "E51"; Value = "  +0.89%  "; ForceText = $false }
)

foreach ($update in $cellUpdates) {
    $range = $ws.Range($update.Cell)
    if ($update.ForceText) {
        # Temporarily mark the cell as Text so the numeric-looking string is not
        # reinterpreted as a number, then restore the default "Normal" style so no
        # stray formatting is left behind on the cell.
        $range.NumberFormat = "@"
        $range.Value = $update.Value
        $range.Style = "Normal"
    } else {
        $range.Value = $update.Value
    }
}
